$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new query columns C and D for rows 2-5 (queries added, values of 25 each)
$rng = $ws.Range("C2:D5")
$rng.Value = 25
$rng.HorizontalAlignment = -4108

# Move active selection to D5
$ws.Range("D5").Select()
